# Auto-generated edits applying the diff to existing_stock sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("existing_stock")

$ws.Cells.Item(11, 5).Value = 0.12617294457189376
$ws.Cells.Item(11, 7).Value = 3162.5000000000005
$ws.Cells.Item(11, 8).Value = 60.500000000000014
$ws.Cells.Item(11, 9).Value = 2.3100000000000005
$ws.Cells.Item(12, 3).Value = "e_w55698557-220"
$ws.Cells.Item(12, 5).Value = 0.11278028565085477
$ws.Cells.Item(13, 3).Value = "e_w758943072-220"
$ws.Cells.Item(13, 5).Value = 0.21146303559535268
$ws.Cells.Item(13, 7).Value = 3162.5
$ws.Cells.Item(14, 3).Value = "e_w240959264-220"
$ws.Cells.Item(15, 3).Value = "e_w281803398-220"
$ws.Cells.Item(15, 5).Value = 0.1047916820839192
$ws.Cells.Item(15, 7).Value = 3162.5000000000005
$ws.Cells.Item(16, 3).Value = "e_w238138373-380"
$ws.Cells.Item(17, 3).Value = "e_CH18-220"
$ws.Cells.Item(17, 5).Value = 0.40037001406053441
$ws.Cells.Item(17, 7).Value = 2749.9999999999995
$ws.Cells.Item(17, 8).Value = 55.000000000000007
$ws.Cells.Item(17, 9).Value = 2.1
$ws.Cells.Item(40, 3).Value = "e_w240959264-220"
$ws.Cells.Item(41, 3).Value = "e_CH18-220"
$ws.Cells.Item(46, 3).Value = "e_w281803398-220"
$ws.Cells.Item(47, 3).Value = "e_w281803398-220"
$ws.Cells.Item(53, 3).Value = "e_w234983117-220"
$ws.Cells.Item(54, 3).Value = "e_w234983117-220"
$ws.Cells.Item(58, 3).Value = "e_CH18-220"
$ws.Cells.Item(59, 3).Value = "e_w240959264-220"
$ws.Cells.Item(69, 3).Value = "e_w238138373-380"
$ws.Cells.Item(101, 3).Value = "e_w234983117-220"
$ws.Cells.Item(122, 3).Value = "e_w238138373-380"
$ws.Cells.Item(124, 3).Value = "e_w758943072-220"
$ws.Cells.Item(125, 3).Value = "e_w238138373-380"
$ws.Cells.Item(140, 3).Value = "e_CH18-220"
$ws.Cells.Item(150, 5).Value = 0.15409499379434963
$ws.Cells.Item(151, 5).Value = 0.15291072157643879
$ws.Cells.Item(152, 5).Value = 0.15521278721895346
$ws.Cells.Item(153, 5).Value = 0.16699472878703805
$ws.Cells.Item(154, 5).Value = 0.17279262926045061
$ws.Cells.Item(155, 5).Value = 0.1586174511333161
$ws.Cells.Item(156, 5).Value = 0.16637932510042519
$ws.Cells.Item(157, 5).Value = 0.21731537653220406
$ws.Cells.Item(158, 5).Value = 0.21471510601685545
$ws.Cells.Item(159, 5).Value = 0.16277902359433066
$ws.Cells.Item(160, 5).Value = 0.16152827258311295
$ws.Cells.Item(161, 5).Value = 0.13606784444360143
$ws.Cells.Item(162, 5).Value = 0.18308354646436523
$ws.Cells.Item(163, 5).Value = 0.15338177131187081
$ws.Cells.Item(164, 5).Value = 0.16554834018408843
$ws.Cells.Item(165, 5).Value = 0.21077006448261207
$ws.Cells.Item(166, 5).Value = 0.17082200788742341
$ws.Cells.Item(167, 5).Value = 0.19614947844032105
$ws.Cells.Item(168, 5).Value = 0.13710016513395351
$ws.Cells.Item(169, 5).Value = 0.20091315882928704
$ws.Cells.Item(170, 5).Value = 0.19828629679661561
$ws.Cells.Item(171, 5).Value = 0.19328994063107527
$ws.Cells.Item(172, 5).Value = 0.1930981018275324
$ws.Cells.Item(173, 5).Value = 0.19782569372870323
$ws.Cells.Item(174, 5).Value = 0.211523174241075
$ws.Cells.Item(175, 3).Value = "elc_spv-CHE_0011"
$ws.Cells.Item(176, 3).Value = "elc_spv-CHE_0023"
$ws.Cells.Item(177, 3).Value = "elc_spv-CHE_0023"
$ws.Cells.Item(178, 3).Value = "elc_spv-CHE_0023"
$ws.Cells.Item(179, 3).Value = "elc_spv-CHE_0023"
$ws.Cells.Item(180, 3).Value = "elc_spv-CHE_0000"
$ws.Cells.Item(181, 3).Value = "elc_spv-CHE_0000"
$ws.Cells.Item(182, 3).Value = "elc_spv-CHE_0000"
$ws.Cells.Item(183, 3).Value = "elc_spv-CHE_0000"
$ws.Cells.Item(184, 3).Value = "elc_spv-CHE_0006"
$ws.Cells.Item(185, 3).Value = "elc_spv-CHE_0006"
$ws.Cells.Item(186, 3).Value = "elc_spv-CHE_0006"
$ws.Cells.Item(187, 3).Value = "elc_spv-CHE_0008"
$ws.Cells.Item(188, 3).Value = "elc_spv-CHE_0008"
$ws.Cells.Item(189, 3).Value = "elc_spv-CHE_0008"
$ws.Cells.Item(190, 3).Value = "elc_spv-CHE_0003"
$ws.Cells.Item(191, 3).Value = "elc_spv-CHE_0003"
$ws.Cells.Item(192, 3).Value = "elc_spv-CHE_0003"
$ws.Cells.Item(193, 3).Value = "elc_spv-CHE_0003"
$ws.Cells.Item(194, 3).Value = "elc_spv-CHE_0003"
$ws.Cells.Item(195, 3).Value = "elc_spv-CHE_0017"
$ws.Cells.Item(196, 3).Value = "elc_spv-CHE_0017"
$ws.Cells.Item(197, 3).Value = "elc_spv-CHE_0017"
$ws.Cells.Item(198, 3).Value = "elc_spv-CHE_0017"
$ws.Cells.Item(199, 3).Value = "elc_spv-CHE_0018"
$ws.Cells.Item(200, 3).Value = "elc_spv-CHE_0018"
$ws.Cells.Item(201, 3).Value = "elc_spv-CHE_0018"
$ws.Cells.Item(202, 3).Value = "elc_spv-CHE_0018"
$ws.Cells.Item(203, 3).Value = "elc_spv-CHE_0018"
$ws.Cells.Item(204, 3).Value = "elc_spv-CHE_0024"
$ws.Cells.Item(205, 3).Value = "elc_spv-CHE_0013"
$ws.Cells.Item(206, 3).Value = "elc_spv-CHE_0012"
$ws.Cells.Item(207, 3).Value = "elc_spv-CHE_0012"
$ws.Cells.Item(208, 3).Value = "elc_spv-CHE_0012"
$ws.Cells.Item(209, 3).Value = "elc_spv-CHE_0012"
$ws.Cells.Item(210, 3).Value = "elc_spv-CHE_0002"
$ws.Cells.Item(211, 3).Value = "elc_spv-CHE_0002"
$ws.Cells.Item(212, 3).Value = "elc_spv-CHE_0025"
